$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Enter Gross Expenditures From 502 Part C" value (D3)
$ws.Range("D3").Value = 122206.33

# Update the "Enter Total Labor Cost From 502 Part L" value (D5)
$ws.Range("D5").Value = 39912.04
